$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo "Miachel" -> "Michael" in the Schumacher question (cell A3)
$ws.Range("A3").Value = "When did Michael Schumacher win his first F1 World Drivers Title?"
